$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8, 9 and 12 have had their species-observation data cyclically
# rotated: row 8 now holds what used to be in row 9, row 9 now holds what
# used to be in row 12, and row 12 now holds what used to be in row 8.
$cols = "A","B","D","E","F","G","H","I","J","P","Q","R"

# Snapshot the current ("before") values first so the writes below don't
# clobber data that still needs to be read for a later column/row.
$row8 = @{}
$row9 = @{}
$row12 = @{}
foreach ($col in $cols) {
    $row8[$col] = $ws.Range($col + "8").Value()
    $row9[$col] = $ws.Range($col + "9").Value()
    $row12[$col] = $ws.Range($col + "12").Value()
}

function Set-RotatedValue($col, $destRow, $value) {
    $target = $ws.Range($col + $destRow)
    if ($value -eq $null -or $value -eq "") {
        $target.Value = ""
    } elseif ($col -eq "I") {
        # The "Antal" column stores numeric-looking counts (e.g. "1") as
        # text, so force a text format before writing to avoid Excel
        # silently re-interpreting the digits as a number.
        $target.NumberFormat = "@"
        $target.Value = $value
    } else {
        $target.Value = $value
    }
}

# row8 <- old row9, row9 <- old row12, row12 <- old row8
foreach ($col in $cols) {
    Set-RotatedValue $col 8 $row9[$col]
    Set-RotatedValue $col 9 $row12[$col]
    Set-RotatedValue $col 12 $row8[$col]
}
